# Remove the "moderate" option from the climate_zone choice list.
#
# The "choices" sheet holds one row per select_one/select_multiple choice.
# The climate_zones list (rows 16-18) currently has: moderate, temperate, hot.
# We delete the "moderate" row (row 16); Excel shifts the remaining rows
# (temperate, hot) up to take rows 16-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")
$ws.Activate()

$ws.Rows.Item(16).Delete()

$ws.Range("A21").Select() | Out-Null
